$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-12 Thursday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-09-13 Friday", 2) | Out-Null
$d.Content.Find.Execute("74×42=", $true, $true, $false, $false, $false, $true, 1, $false, "31×76=", 2) | Out-Null
$d.Content.Find.Execute("25×21=", $true, $true, $false, $false, $false, $true, 1, $false, "69×66=", 2) | Out-Null
$d.Content.Find.Execute("17×28=", $true, $true, $false, $false, $false, $true, 1, $false, "57×17=", 2) | Out-Null
$d.Content.Find.Execute("40×37=", $true, $true, $false, $false, $false, $true, 1, $false, "50×25=", 2) | Out-Null
$d.Content.Find.Execute("13×11=", $true, $true, $false, $false, $false, $true, 1, $false, "60×19=", 2) | Out-Null
$d.Content.Find.Execute("38×80=", $true, $true, $false, $false, $false, $true, 1, $false, "47×45=", 2) | Out-Null
$d.Content.Find.Execute("37×99=", $true, $true, $false, $false, $false, $true, 1, $false, "15×17=", 2) | Out-Null
$d.Content.Find.Execute("13×14=", $true, $true, $false, $false, $false, $true, 1, $false, "32×72=", 2) | Out-Null
$d.Content.Find.Execute("84×58=", $true, $true, $false, $false, $false, $true, 1, $false, "31×94=", 2) | Out-Null
$d.Content.Find.Execute("14×30=", $true, $true, $false, $false, $false, $true, 1, $false, "17×64=", 2) | Out-Null
$d.Content.Find.Execute("38×20=", $true, $true, $false, $false, $false, $true, 1, $false, "35×99=", 2) | Out-Null
$d.Content.Find.Execute("90×46=", $true, $true, $false, $false, $false, $true, 1, $false, "16×70=", 2) | Out-Null
$d.Content.Find.Execute("77×27=", $true, $true, $false, $false, $false, $true, 1, $false, "25×73=", 2) | Out-Null
$d.Content.Find.Execute("26×80=", $true, $true, $false, $false, $false, $true, 1, $false, "83×81=", 2) | Out-Null
$d.Content.Find.Execute("51×25=", $true, $true, $false, $false, $false, $true, 1, $false, "39×89=", 2) | Out-Null
$d.Content.Find.Execute("74×92=", $true, $true, $false, $false, $false, $true, 1, $false, "54×86=", 2) | Out-Null
$d.Content.Find.Execute("93×66=", $true, $true, $false, $false, $false, $true, 1, $false, "51×95=", 2) | Out-Null
$d.Content.Find.Execute("28×95=", $true, $true, $false, $false, $false, $true, 1, $false, "39×68=", 2) | Out-Null
$d.Content.Find.Execute("72×21=", $true, $true, $false, $false, $false, $true, 1, $false, "89×32=", 2) | Out-Null
$d.Content.Find.Execute("87×72=", $true, $true, $false, $false, $false, $true, 1, $false, "97×13=", 2) | Out-Null
$d.Content.Find.Execute("59×58=", $true, $true, $false, $false, $false, $true, 1, $false, "99×37=", 2) | Out-Null
$d.Content.Find.Execute("18×61=", $true, $true, $false, $false, $false, $true, 1, $false, "55×94=", 2) | Out-Null
$d.Content.Find.Execute("76×72=", $true, $true, $false, $false, $false, $true, 1, $false, "88×74=", 2) | Out-Null
$d.Content.Find.Execute("44×12=", $true, $true, $false, $false, $false, $true, 1, $false, "33×60=", 2) | Out-Null
$d.Content.Find.Execute("24×20=", $true, $true, $false, $false, $false, $true, 1, $false, "82×78=", 2) | Out-Null
